# Update "想去人数" (column F) counts for several events on both the
# "展览" sheet and the "全部类型" sheet, matching the regenerated data
# output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4356
$wsExhibit.Range("F3").Value = 2451
$wsExhibit.Range("F6").Value = 48
$wsExhibit.Range("F12").Value = 1624
$wsExhibit.Range("F14").Value = 3425

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4356
$wsAll.Range("F3").Value = 2451
$wsAll.Range("F7").Value = 48
$wsAll.Range("F16").Value = 1624
$wsAll.Range("F18").Value = 3425
